$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 842.8570999999999
$ws.Range("I12").Value = 620.25
$ws.Range("J12").Value = 1139.6666
$ws.Range("K12").Value = 620.25
$ws.Range("L12").Value = 1139.6666
$ws.Range("M12").Value = -450.25
$ws.Range("N12").Value = -1479.6666

$ws.Range("H80").Value = 4312495.5
$ws.Range("J80").Value = 2362.9167
$ws.Range("L80").Value = 7088.750100000001
$ws.Range("N80").Value = -9084.750100000001

$ws.Range("H83").Value = 4312495.5
$ws.Range("J83").Value = 2362.9167
$ws.Range("L83").Value = 21266.2503
$ws.Range("N83").Value = -31250.2503

$ws.Range("H103").Value = 449.5
$ws.Range("I103").Value = 450
$ws.Range("J103").Value = 449
$ws.Range("K103").Value = 1350
$ws.Range("L103").Value = 1347
$ws.Range("M103").Value = -764
$ws.Range("N103").Value = -2519

$ws.Range("H132").Value = 35812.285
$ws.Range("I132").Value = 41771.156
$ws.Range("J132").Value = 7011.0835
$ws.Range("K132").Value = 125313.468
$ws.Range("L132").Value = 21033.2505
$ws.Range("M132").Value = -122783.468
$ws.Range("N132").Value = -26093.2505

$ws.Range("H138").Value = 7211.3193
$ws.Range("I138").Value = 14206.277
$ws.Range("J138").Value = 2869.6206
$ws.Range("K138").Value = 42618.831
$ws.Range("L138").Value = 8608.861800000001
$ws.Range("M138").Value = -37478.831
$ws.Range("N138").Value = -18888.8618

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1927.4166
$ws.Range("I122").Value = 1213
$ws.Range("K122").Value = 3639
$ws.Range("M122").Value = -1189

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8744.632
$ws.Range("I99").Value = 4075.9333
$ws.Range("K99").Value = 4075.9333
$ws.Range("M99").Value = -2577.9333

$ws.Range("H103").Value = 72181.39999999999
$ws.Range("J103").Value = 72181.39999999999
$ws.Range("L103").Value = 72181.39999999999
$ws.Range("N103").Value = -74525.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 142082.77
$ws.Range("I31").Value = 348633.12
$ws.Range("K31").Value = 348633.12
$ws.Range("M31").Value = -348338.12

$ws.Range("H34").Value = 142082.77
$ws.Range("I34").Value = 348633.12
$ws.Range("K34").Value = 348633.12
$ws.Range("M34").Value = -348431.12

$ws.Range("H43").Value = 97612
$ws.Range("J43").Value = 97612
$ws.Range("L43").Value = 97612
$ws.Range("N43").Value = -97980

$ws.Range("H62").Value = 4834.8335
$ws.Range("J62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248

$ws.Range("H65").Value = 4834.8335
$ws.Range("J65").Value = 7000
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240

$ws.Range("H86").Value = 207786.4
$ws.Range("I86").Value = 7993.8335
$ws.Range("K86").Value = 7993.8335
$ws.Range("M86").Value = -6870.8335

$ws.Range("H89").Value = 207786.4
$ws.Range("I89").Value = 7993.8335
$ws.Range("K89").Value = 39969.1675
$ws.Range("M89").Value = -34353.1675

$ws.Range("H101").Value = 97612
$ws.Range("J101").Value = 97612
$ws.Range("L101").Value = 97612
$ws.Range("N101").Value = -104102

$ws.Range("H122").Value = 3065.1667
$ws.Range("I122").Value = 1972.75
$ws.Range("K122").Value = 5918.25
$ws.Range("M122").Value = -3468.25

$ws.Range("H132").Value = 50012290
$ws.Range("I132").Value = 62514516
$ws.Range("J132").Value = 19237574
$ws.Range("K132").Value = 187543548
$ws.Range("L132").Value = 57712722
$ws.Range("M132").Value = -187541018
$ws.Range("N132").Value = -57717782

$ws.Range("H134").Value = 15538.88
$ws.Range("I134").Value = 19203.316
$ws.Range("K134").Value = 57609.948
$ws.Range("M134").Value = -55074.948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 31068.166
$ws.Range("I3").Value = 29602.25
$ws.Range("K3").Value = 88806.75
$ws.Range("M3").Value = -88694.75

$ws.Range("H68").Value = 10109.23
$ws.Range("I68").Value = 2224.75
$ws.Range("K68").Value = 6674.25
$ws.Range("M68").Value = -5863.25

$ws.Range("H71").Value = 10109.23
$ws.Range("I71").Value = 2224.75
$ws.Range("K71").Value = 20022.75
$ws.Range("M71").Value = -15966.75

$ws.Range("H87").Value = 16502.8
$ws.Range("I87").Value = 10827.333
$ws.Range("K87").Value = 32481.999
$ws.Range("M87").Value = -31233.999

$ws.Range("H90").Value = 16502.8
$ws.Range("I90").Value = 10827.333
$ws.Range("K90").Value = 97445.997
$ws.Range("M90").Value = -91205.997

$ws.Range("H107").Value = 208.4
$ws.Range("J107").Value = 213
$ws.Range("L107").Value = 639
$ws.Range("N107").Value = -4479

$ws.Range("H113").Value = 1269.871
$ws.Range("I113").Value = 466.66666
$ws.Range("J113").Value = 1355.9286
$ws.Range("K113").Value = 1399.99998
$ws.Range("L113").Value = 4067.7858
$ws.Range("M113").Value = 770.0000199999999
$ws.Range("N113").Value = -8407.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 78536.8
$ws.Range("J105").Value = 78536.8
$ws.Range("L105").Value = 78536.8
$ws.Range("N105").Value = -85524.8

$ws.Range("H126").Value = 1193759.1
$ws.Range("I126").Value = 2085678.5
$ws.Range("J126").Value = 4533.3335
$ws.Range("K126").Value = 6257035.5
$ws.Range("L126").Value = 13600.0005
$ws.Range("M126").Value = -6254565.5
$ws.Range("N126").Value = -18540.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 450.72
$ws.Range("I16").Value = 460.7143
$ws.Range("J16").Value = 438
$ws.Range("K16").Value = 460.7143
$ws.Range("L16").Value = 438
$ws.Range("M16").Value = -290.7143
$ws.Range("N16").Value = -778

$ws.Range("H40").Value = 2487.4614
$ws.Range("I40").Value = 2417
$ws.Range("J40").Value = 2875
$ws.Range("K40").Value = 2417
$ws.Range("L40").Value = 2875
$ws.Range("M40").Value = -2281
$ws.Range("N40").Value = -3147

$ws.Range("H46").Value = 903
$ws.Range("I46").Value = 704
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 704
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -516
$ws.Range("N46").Value = -1876

$ws.Range("H136").Value = 40706.824
$ws.Range("J136").Value = 89263.2
$ws.Range("L136").Value = 267789.6
$ws.Range("N136").Value = -272889.6

$ws.Range("H140").Value = 85426
$ws.Range("J140").Value = 85426
$ws.Range("L140").Value = 85426
$ws.Range("N140").Value = -95786

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1896.5769
$ws.Range("I107").Value = 1084.238
$ws.Range("K107").Value = 3252.714
$ws.Range("M107").Value = -1332.714

$ws.Range("H132").Value = 7297997.5
$ws.Range("I132").Value = 9555068
$ws.Range("K132").Value = 28665204
$ws.Range("M132").Value = -28662674
